# "Generate Report for Archive"
# The localization-status report is regenerated; the row that used to be in
# slot 5 (ceaa8c31-441d-428a-91f5-de4fcc64abec) now reports as slot 3, and
# the two rows that used to be in slots 3/4 (dc4b74d3.../85049ea7...) shift
# down to slots 4/5 on every sheet (Overview, zh-cn, de-de).
#
# We apply this as a plain 3-way rotation of row 3 -> 4 -> 5 -> 3 restricted
# to the columns that actually carry per-file data (the row-number / style
# columns and any column that is identical across the three rows, such as
# the blank "Source Path" column, are left untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A, B, E, F, G vary per file (C, D are shared)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$cols = @("A", "B", "E", "F", "G")
$before = @{}
foreach ($c in $cols) {
    $before[$c + "3"] = $ws.Range($c + "3").Value2
    $before[$c + "4"] = $ws.Range($c + "4").Value2
    $before[$c + "5"] = $ws.Range($c + "5").Value2
}
foreach ($c in $cols) {
    $ws.Range($c + "3").Value2 = $before[$c + "5"]
    $ws.Range($c + "4").Value2 = $before[$c + "3"]
    $ws.Range($c + "5").Value2 = $before[$c + "4"]
}

# Hyperlinks in column B carry a "display" label independent of the cell
# text; rotate those the same way while keeping each hyperlink's actual
# target address (and therefore its relationship id) pinned to its row.
$linkAddr = @{}
$linkDisp = @{}
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $h = $ws.Hyperlinks.Item($i)
    $ref = $h.Range.Address($false, $false)
    $linkAddr[$ref] = $h.Address
    $linkDisp[$ref] = $h.TextToDisplay
}
$d3 = $linkDisp["B3"]
$d4 = $linkDisp["B4"]
$d5 = $linkDisp["B5"]
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $linkAddr["B2"], "", "", $linkDisp["B2"])
$ws.Hyperlinks.Add($ws.Range("B3"), $linkAddr["B3"], "", "", $d5)
$ws.Hyperlinks.Add($ws.Range("B4"), $linkAddr["B4"], "", "", $d3)
$ws.Hyperlinks.Add($ws.Range("B5"), $linkAddr["B5"], "", "", $d4)

# ---------------------------------------------------------------------
# Sheets "zh-cn" and "de-de": columns A, C, F, G vary per file
# ---------------------------------------------------------------------
$langSheets = @("zh-cn", "de-de")
foreach ($sheetName in $langSheets) {
    $ws2 = $wb.Worksheets.Item($sheetName)

    $cols2 = @("A", "C", "F", "G")
    $before2 = @{}
    foreach ($c in $cols2) {
        $before2[$c + "3"] = $ws2.Range($c + "3").Value2
        $before2[$c + "4"] = $ws2.Range($c + "4").Value2
        $before2[$c + "5"] = $ws2.Range($c + "5").Value2
    }
    foreach ($c in $cols2) {
        $ws2.Range($c + "3").Value2 = $before2[$c + "5"]
        $ws2.Range($c + "4").Value2 = $before2[$c + "3"]
        $ws2.Range($c + "5").Value2 = $before2[$c + "4"]
    }

    # Hyperlinks live in column A on these sheets.
    $linkAddr2 = @{}
    $linkDisp2 = @{}
    for ($i = 1; $i -le $ws2.Hyperlinks.Count; $i++) {
        $h2 = $ws2.Hyperlinks.Item($i)
        $ref2 = $h2.Range.Address($false, $false)
        $linkAddr2[$ref2] = $h2.Address
        $linkDisp2[$ref2] = $h2.TextToDisplay
    }
    $a3 = $linkDisp2["A3"]
    $a4 = $linkDisp2["A4"]
    $a5 = $linkDisp2["A5"]
    $ws2.Hyperlinks.Delete()
    $ws2.Hyperlinks.Add($ws2.Range("A2"), $linkAddr2["A2"], "", "", $linkDisp2["A2"])
    $ws2.Hyperlinks.Add($ws2.Range("A3"), $linkAddr2["A3"], "", "", $a5)
    $ws2.Hyperlinks.Add($ws2.Range("A4"), $linkAddr2["A4"], "", "", $a3)
    $ws2.Hyperlinks.Add($ws2.Range("A5"), $linkAddr2["A5"], "", "", $a4)
}
